# issue #5: stock data from json to db
#
# Adds three "pipeline" metadata columns (category, source_file, index) to
# the 股票 (stock) sheet, matching the shape already produced for the other
# property sheets by the JSON->DB export: a "category" column right after
# "property_category", and "source_file"/"index" columns appended after
# "legislator_id".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- insert the new "category" column between "property_category" (H) and
#     "date" (old I) -------------------------------------------------------
$ws.Columns("I:I").Insert()
$ws.Range("I1").Value = "category"
$ws.Range("I2").Value = "normal"

# --- append "source_file" and "index" columns after "legislator_id" (L) --
$ws.Range("H1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("N1").PasteSpecial(-4122)

$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

$ws.Range("M2").Value = "tmp63271"
$ws.Range("N2").Value = 94
